# repull data, push all data, mean calculation
# Update the "dSF" (F) column values for several rows to reflect repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -8
$ws.Range("F4").Value = -6
$ws.Range("F7").Value = -10
$ws.Range("F19").Value = -4
$ws.Range("F33").Value = -3
$ws.Range("F34").Value = 5
$ws.Range("F35").Value = 6
$ws.Range("F36").Value = -5
$ws.Range("F39").Value = -5
$ws.Range("F45").Value = -1
